# Daily Report update: 2026-02-13
# Appends the next day's 22-row depository block (11 companies x
# Registered/Eligible) to the "Daily_Data" sheet, duplicating the most
# recent block (rows 596:617, serial date 46064) and bumping the date
# serial to 46065 for the new block (rows 618:639).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

# Find the last used row of the existing data block.
$lastRow = $ws.UsedRange.Rows.Count

# The most recent day occupies a 22-row block ending at $lastRow.
$blockSize = 22
$srcFirst = $lastRow - $blockSize + 1
$srcLast = $lastRow

$dstFirst = $lastRow + 1
$dstLast = $lastRow + $blockSize

# Duplicate the prior day's block (columns A:H) into the new rows.
$src = $ws.Range("A" + $srcFirst + ":H" + $srcLast)
$dst = $ws.Range("A" + $dstFirst)
$src.Copy($dst)

# New reporting date serial (2026-02-12 underlying serial, reported as the
# "2026-02-13" Daily Report): one day after the previous block's date.
$prevDate = $ws.Range("A" + $srcFirst).Value2
$newDate = $prevDate + 1
$ws.Range("A" + $dstFirst + ":A" + $dstLast).Value = $newDate
